# Update cryptos list values per upstream diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.918.07"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "2.776.72"
$ws.Range("E3").Value = "  -1.65%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "356.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.36%  "
$ws.Range("E7").Value = "  +1.79%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.588"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.89%  "
$ws.Range("E10").Value = "  -4.54%  "
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("E13").Value = "  -3.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("D15").Value = "3.214.63"
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("D16").Value = "2.775.11"
$ws.Range("E16").Value = "  -2.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.934"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.11%  "
$ws.Range("D18").Value = "51.826.56"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.27%  "
$ws.Range("E21").Value = "  -3.78%  "
$ws.Range("D22").Value = "0.0₃0973"
$ws.Range("E22").Value = "  -2.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "274.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  -2.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.49"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.23%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.10%  "
$ws.Range("E29").Value = "  -1.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.144"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0466"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.12%  "
$ws.Range("E32").Value = "  +1.12%  "
$ws.Range("E33").Value = "  -0.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.69%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0839"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.99%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.29"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("E38").Value = "  -0.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.28%  "
$ws.Range("E40").Value = "  -4.99%  "
$ws.Range("E41").Value = "  -1.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.114"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "121.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.52%  "
$ws.Range("D46").Value = "2.049.38"
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("E47").Value = "  -3.49%  "
$ws.Range("E48").Value = "  -1.99%  "
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.920"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.26%  "
$ws.Range("E51").Value = "  -0.32%  "
